# Apply the recorded change set to the "Artfynd" worksheet.
# Rows 14-16 effectively rotate their species-observation data
# (row14 <- old row15, row15 <- old row16, row16 <- old row14),
# with a couple of individually-tweaked values (B15, AY14/AY15/AY16).
# Row 13 only has a small numeric tweak to B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 ---
$ws.Range("B13").Value = 93334

# --- Row 14 ---
$ws.Range("A14").Value = 112231491
$ws.Range("B14").Value = 56575
$ws.Range("E14").Value = 103021
$ws.Range("F14").Value = "Talltita"
$ws.Range("G14").Value = "Poecile montanus"
$ws.Range("H14").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q14").Value = 478579
$ws.Range("R14").Value = 6556322
$ws.Range("S14").Value = 10
$ws.Range("AY14").Value = "Länsstyrelsen i Örebro län, inventering"

# --- Row 15 ---
$ws.Range("A15").Value = 112205187
$ws.Range("B15").Value = 90810
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 4363
$ws.Range("F15").Value = "Zontaggsvamp"
$ws.Range("G15").Value = "Hydnellum concrescens"
$ws.Range("H15").Value = "(Pers.) Banker"
$ws.Range("J15").Value = "mycel"
# K15 already holds an (empty-text) value identical to the target - leave untouched.
$ws.Range("L15").Value = $null
$ws.Range("M15").Value = $null
$ws.Range("Q15").Value = 478586
$ws.Range("R15").Value = 6556137
# Leading "'" forces a literal (non-formula) empty-text value instead of clearing the cell;
# reset the style afterwards so the implicit quote-prefix formatting doesn't stick.
$ws.Range("AF15").Value = "'"
$ws.Range("AF15").Style = "Normal"
$ws.Range("AH15").Value = "Blåbärsbarrskog"
$ws.Range("AI15").Value = "i yta bökad av vildsvin"

# --- Row 16 ---
$ws.Range("A16").Value = 112269209
$ws.Range("B16").Value = 56446
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 100049
$ws.Range("F16").Value = "Spillkråka"
$ws.Range("G16").Value = "Dryocopus martius"
$ws.Range("H16").Value = "(Linnaeus, 1758)"
$ws.Range("J16").Value = $null
$ws.Range("L16").Value = "'"
$ws.Range("L16").Style = "Normal"
$ws.Range("M16").Value = "lockläte, övriga läten"
$ws.Range("Q16").Value = 478539
$ws.Range("R16").Value = 6556219
$ws.Range("S16").Value = 100
$ws.Range("AF16").Value = $null
$ws.Range("AH16").Value = $null
$ws.Range("AI16").Value = $null
$ws.Range("AY16").Value = "'"
$ws.Range("AY16").Style = "Normal"
